$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.241.54'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +3.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.778.67'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.31%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '339.00'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3827'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -3.39%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3429'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.86'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.149'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07404'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '23.63'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.49%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.9993'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.436'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.314'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.82%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.780.47'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.12%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001081'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.39%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06682'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '82.60'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.004'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.44'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.420'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.181.23'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.20%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.386'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.76'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.430'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.411'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.58%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '153.76'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.978.86'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.36%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '134.59'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.34%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.027'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.089'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08842'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '12.78'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02415'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.6841'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.337'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06378'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.32%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2161'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.246'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.506'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.04%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.301'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.93%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.002'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.09'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6277'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.87%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.858'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.63%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.075'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07501'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.99%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.194'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.43%  '
